$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B65: change from text "3" to numeric 3
$ws.Range("B65").Value = 3

# Add new row 66
$ws.Range("A66").Value = "Ying Tang"
$ws.Range("B66").NumberFormat = "@"
$ws.Range("B66").Value = "2"
$ws.Range("B66").Style = "Normal"
$ws.Range("C66").Value = " needs work,some missing pieces"
$ws.Range("D66").Value = "CRT"
$ws.Range("E66").Value = "WRI"
$ws.Range("F66").Value = "c1f109d6-e04a-469c-a254-426c0826b7a8"
$ws.Range("G66").Value = "BJcAWaeCW_annotated.xlsx"
$ws.Range("H66").Value = "The paper organization needs work; there are also some missing pieces to put the NN training together."
